# Fruta / hortaliza, semanal
# Insert a new week of "Chirimoya" price data (3 rows) ahead of the existing
# data block, shifting the previously-recorded rows down by three positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows just above row 151; this pushes the existing rows
# 151-199 down to 154-202 (carrying their formatting, including the date
# number format on column D) and extends the sheet's used range / dimension
# to A1:T202 automatically.
$ws.Rows("151:153").Insert()

# Populate the 3 newly-inserted rows with this week's data. All of the
# non-varying descriptive columns (A, B, C, E-J) mirror the rest of the
# "Chirimoya" / "Mercado Mayorista Lo Valledor de Santiago" block.

# Row 151: Cultivar IV Región - Especial - Provincia de Limarí
$ws.Range("A151").Value = 6
$ws.Range("B151").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C151").Value = "Metropolitana"
$ws.Range("D151").Value = 44524
$ws.Range("E151").Value = 13
$ws.Range("F151").Value = "Fruta"
$ws.Range("G151").Value = 100107
$ws.Range("H151").Value = "Otros"
$ws.Range("I151").Value = 100107002
$ws.Range("J151").Value = "Chirimoya"
$ws.Range("K151").Value = "Cultivar IV Región"
$ws.Range("L151").Value = "Especial"
$ws.Range("M151").Value = 170
$ws.Range("N151").Value = 2300
$ws.Range("O151").Value = 2300
$ws.Range("P151").Value = 2300
$ws.Range("Q151").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R151").Value = "Provincia de Limarí"
$ws.Range("S151").Value = 2300
$ws.Range("T151").Value = 1

# Row 152: Cultivar IV Región - Primera - Provincia de Limarí
$ws.Range("A152").Value = 6
$ws.Range("B152").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C152").Value = "Metropolitana"
$ws.Range("D152").Value = 44524
$ws.Range("E152").Value = 13
$ws.Range("F152").Value = "Fruta"
$ws.Range("G152").Value = 100107
$ws.Range("H152").Value = "Otros"
$ws.Range("I152").Value = 100107002
$ws.Range("J152").Value = "Chirimoya"
$ws.Range("K152").Value = "Cultivar IV Región"
$ws.Range("L152").Value = "Primera"
$ws.Range("M152").Value = 230
$ws.Range("N152").Value = 2000
$ws.Range("O152").Value = 2000
$ws.Range("P152").Value = 2000
$ws.Range("Q152").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R152").Value = "Provincia de Limarí"
$ws.Range("S152").Value = 2000
$ws.Range("T152").Value = 1

# Row 153: Cultivar IV Región - Segunda - Provincia de Limarí
$ws.Range("A153").Value = 6
$ws.Range("B153").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C153").Value = "Metropolitana"
$ws.Range("D153").Value = 44524
$ws.Range("E153").Value = 13
$ws.Range("F153").Value = "Fruta"
$ws.Range("G153").Value = 100107
$ws.Range("H153").Value = "Otros"
$ws.Range("I153").Value = 100107002
$ws.Range("J153").Value = "Chirimoya"
$ws.Range("K153").Value = "Cultivar IV Región"
$ws.Range("L153").Value = "Segunda"
$ws.Range("M153").Value = 200
$ws.Range("N153").Value = 1600
$ws.Range("O153").Value = 1600
$ws.Range("P153").Value = 1600
$ws.Range("Q153").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R153").Value = "Provincia de Limarí"
$ws.Range("S153").Value = 1600
$ws.Range("T153").Value = 1
